# Update Leve profit-calculation sheets with refreshed market-board prices.
# Each sheet has the same layout: columns H-N hold price/profit figures
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 950.3158
$ws.Range("J9").Value = 854
$ws.Range("L9").Value = 854
$ws.Range("N9").Value = -1192
$ws.Range("H17").Value = 1998.091
$ws.Range("J17").Value = 3094.1428
$ws.Range("L17").Value = 9282.428400000001
$ws.Range("N17").Value = -9618.428400000001
$ws.Range("H33").Value = 1207.5883
$ws.Range("I33").Value = 219.16667
$ws.Range("K33").Value = 219.16667
$ws.Range("M33").Value = 9.833329999999989
$ws.Range("H76").Value = 9697.6
$ws.Range("I76").Value = 9622
$ws.Range("J76").Value = 10000
$ws.Range("K76").Value = 9622
$ws.Range("L76").Value = 10000
$ws.Range("M76").Value = -9307
$ws.Range("N76").Value = -10630
$ws.Range("H79").Value = 9697.6
$ws.Range("I79").Value = 9622
$ws.Range("J79").Value = 10000
$ws.Range("K79").Value = 9622
$ws.Range("L79").Value = 10000
$ws.Range("M79").Value = -8530
$ws.Range("N79").Value = -12184
$ws.Range("H92").Value = 23693.375
$ws.Range("I92").Value = 930.13794
$ws.Range("K92").Value = 930.13794
$ws.Range("M92").Value = 317.86206
$ws.Range("H96").Value = 1614.2778
$ws.Range("I96").Value = 848.2222
$ws.Range("J96").Value = 2380.3333
$ws.Range("K96").Value = 2544.6666
$ws.Range("L96").Value = 7140.999899999999
$ws.Range("M96").Value = -1171.6666
$ws.Range("N96").Value = -9886.999899999999
$ws.Range("H99").Value = 1028.5714
$ws.Range("J99").Value = 1340.3334
$ws.Range("L99").Value = 4021.0002
$ws.Range("N99").Value = -7017.0002
$ws.Range("H106").Value = 3774.4119
$ws.Range("I106").Value = 3873.7693
$ws.Range("K106").Value = 3873.7693
$ws.Range("M106").Value = -3242.7693
$ws.Range("H107").Value = 6082.8335
$ws.Range("I107").Value = 5912.1875
$ws.Range("J107").Value = 7448
$ws.Range("K107").Value = 5912.1875
$ws.Range("L107").Value = 7448
$ws.Range("M107").Value = -3992.1875
$ws.Range("N107").Value = -11288
$ws.Range("H118").Value = 1621.3334
$ws.Range("I118").Value = 1266.5
$ws.Range("J118").Value = 1905.2
$ws.Range("K118").Value = 3799.5
$ws.Range("L118").Value = 5715.6
$ws.Range("M118").Value = -2142.5
$ws.Range("N118").Value = -9029.6
$ws.Range("H131").Value = 8110.467
$ws.Range("I131").Value = 3043.2856
$ws.Range("J131").Value = 19933.889
$ws.Range("K131").Value = 9129.856800000001
$ws.Range("L131").Value = 59801.667
$ws.Range("M131").Value = -4089.856800000001
$ws.Range("N131").Value = -69881.667
$ws.Range("H132").Value = 2526.8
$ws.Range("I132").Value = 2544.1135
$ws.Range("K132").Value = 7632.3405
$ws.Range("M132").Value = -5102.3405

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2310.2173
$ws.Range("I97").Value = 772.8182
$ws.Range("J97").Value = 3719.5
$ws.Range("K97").Value = 772.8182
$ws.Range("L97").Value = 3719.5
$ws.Range("M97").Value = -276.8182
$ws.Range("N97").Value = -4711.5
$ws.Range("H122").Value = 2414.92
$ws.Range("I122").Value = 2335.2727
$ws.Range("K122").Value = 7005.8181
$ws.Range("M122").Value = -4555.8181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 12827.105
$ws.Range("I107").Value = 15443.6
$ws.Range("J107").Value = 3015.25
$ws.Range("K107").Value = 15443.6
$ws.Range("L107").Value = 3015.25
$ws.Range("M107").Value = -13523.6
$ws.Range("N107").Value = -6855.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 20000
$ws.Range("J48").Value = 20000
$ws.Range("L48").Value = 20000
$ws.Range("N48").Value = -20952
$ws.Range("H58").Value = 1958.2858
$ws.Range("I58").Value = 741.8
$ws.Range("J58").Value = 4999.5
$ws.Range("K58").Value = 741.8
$ws.Range("L58").Value = 4999.5
$ws.Range("M58").Value = -538.8
$ws.Range("N58").Value = -5405.5
$ws.Range("H70").Value = 70000
$ws.Range("J70").Value = 70000
$ws.Range("L70").Value = 70000
$ws.Range("N70").Value = -70630
$ws.Range("H73").Value = 70000
$ws.Range("J73").Value = 70000
$ws.Range("L73").Value = 70000
$ws.Range("N73").Value = -72184
$ws.Range("H105").Value = 1769.5555
$ws.Range("J105").Value = 1599.6666
$ws.Range("L105").Value = 1599.6666
$ws.Range("N105").Value = -5093.6666
$ws.Range("H136").Value = 1958.2858
$ws.Range("I136").Value = 741.8
$ws.Range("J136").Value = 4999.5
$ws.Range("K136").Value = 2225.4
$ws.Range("L136").Value = 14998.5
$ws.Range("M136").Value = 324.6000000000004
$ws.Range("N136").Value = -20098.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 18135.818
$ws.Range("I56").Value = 18135.818
$ws.Range("K56").Value = 18135.818
$ws.Range("M56").Value = -17605.818
$ws.Range("H131").Value = 1795.0952
$ws.Range("J131").Value = 2059.9795
$ws.Range("L131").Value = 6179.9385
$ws.Range("N131").Value = -16259.9385
$ws.Range("H133").Value = 1000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 1000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 3000
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -13120
$ws.Range("H134").Value = 6250.154
$ws.Range("I134").Value = 6250.154
$ws.Range("K134").Value = 18750.462
$ws.Range("M134").Value = -13680.462

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 43042
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H50").Value = 43042
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H102").Value = 4499.205
$ws.Range("I102").Value = 4133.472
$ws.Range("K102").Value = 4133.472
$ws.Range("M102").Value = -2511.472
$ws.Range("H122").Value = 46327.32
$ws.Range("I122").Value = 64844.176
$ws.Range("J122").Value = 6979
$ws.Range("K122").Value = 194532.528
$ws.Range("L122").Value = 20937
$ws.Range("M122").Value = -192082.528
$ws.Range("N122").Value = -25837
$ws.Range("H126").Value = 3278
$ws.Range("I126").Value = 3278
$ws.Range("K126").Value = 9834
$ws.Range("M126").Value = -7364
$ws.Range("H132").Value = 3751.8462
$ws.Range("I132").Value = 3643.257
$ws.Range("J132").Value = 4702
$ws.Range("K132").Value = 10929.771
$ws.Range("L132").Value = 14106
$ws.Range("M132").Value = -8399.771000000001
$ws.Range("N132").Value = -19166
$ws.Range("H135").Value = 90375
$ws.Range("J135").Value = 90375
$ws.Range("L135").Value = 90375
$ws.Range("N135").Value = -100515

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4503.857
$ws.Range("I7").Value = 4520.4
$ws.Range("K7").Value = 4520.4
$ws.Range("M7").Value = -4408.4
$ws.Range("H55").Value = 1316.5555
$ws.Range("I55").Value = 270.7647
$ws.Range("J55").Value = 3094.4
$ws.Range("K55").Value = 270.7647
$ws.Range("L55").Value = 3094.4
$ws.Range("M55").Value = -97.7647
$ws.Range("N55").Value = -3440.4
$ws.Range("H126").Value = 4503.857
$ws.Range("I126").Value = 4520.4
$ws.Range("K126").Value = 13561.2
$ws.Range("M126").Value = -11091.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4764872
$ws.Range("I81").Value = 6496476.5
$ws.Range("J81").Value = 2960.25
$ws.Range("K81").Value = 12992953
$ws.Range("L81").Value = 5920.5
$ws.Range("M81").Value = -12991892
$ws.Range("N81").Value = -8042.5
$ws.Range("H84").Value = 4764872
$ws.Range("I84").Value = 6496476.5
$ws.Range("J84").Value = 2960.25
$ws.Range("K84").Value = 64964765
$ws.Range("L84").Value = 29602.5
$ws.Range("M84").Value = -64959461
$ws.Range("N84").Value = -40210.5
$ws.Range("H122").Value = 6733.7344
$ws.Range("I122").Value = 6782.5
$ws.Range("J122").Value = 6640.636
$ws.Range("K122").Value = 20347.5
$ws.Range("L122").Value = 19921.908
$ws.Range("M122").Value = -17897.5
$ws.Range("N122").Value = -24821.908
$ws.Range("H132").Value = 1311.75
$ws.Range("I132").Value = 1000.6667
$ws.Range("K132").Value = 3002.0001
